{"js": "// This script applies the set of literal text replacements that correspond\n// to the day's updated date heading and the newly generated division\n// answers in the table. Every `<w:t>` run in the document is replaced\n// exactly once, in document order, using the Word JS API's body.search().\n\nconst replacements = [\n  [\"2025-09-11 Thursday\", \"2025-09-12 Friday\"],\n  [\"79\u00f78=9, 7\", \"94\u00f79=10, 4\"],\n  [\"49\u00f73=16, 1\", \"74\u00f74=18, 2\"],\n  [\"86\u00f79=9, 5\", \"58\u00f77=8, 2\"],\n  [\"20\u00f78=2, 4\", \"75\u00f77=10, 5\"],\n  [\"52\u00f77=7, 3\", \"70\u00f73=23, 1\"],\n  [\"65\u00f79=7, 2\", \"51\u00f74=12, 3\"],\n  [\"41\u00f76=6, 5\", \"76\u00f73=25, 1\"],\n  [\"66\u00f74=16, 2\", \"22\u00f77=3, 1\"],\n  [\"92\u00f73=30, 2\", \"41\u00f76=6, 5\"],\n  [\"33\u00f77=4, 5\", \"87\u00f75=17, 2\"],\n  [\"75\u00f79=8, 3\", \"43\u00f74=10, 3\"],\n  [\"68\u00f74=17, 0\", \"32\u00f73=10, 2\"],\n  [\"27\u00f77=3, 6\", \"22\u00f76=3, 4\"],\n  [\"15\u00f76=2, 3\", \"16\u00f74=4, 0\"],\n  [\"86\u00f72=43, 0\", \"43\u00f77=6, 1\"],\n  [\"82\u00f73=27, 1\", \"33\u00f74=8, 1\"],\n  [\"65\u00f72=32, 1\", \"23\u00f77=3, 2\"],\n  [\"90\u00f75=18, 0\", \"89\u00f74=22, 1\"],\n  [\"29\u00f76=4, 5\", \"43\u00f72=21, 1\"],\n  [\"34\u00f72=17, 0\", \"65\u00f74=16, 1\"],\n  [\"74\u00f73=24, 2\", \"77\u00f75=15, 2\"],\n  [\"80\u00f72=40, 0\", \"16\u00f75=3, 1\"],\n  [\"13\u00f75=2, 3\", \"31\u00f78=3, 7\"],\n  [\"26\u00f77=3, 5\", \"31\u00f72=15, 1\"],\n  [\"96\u00f77=13, 5\", \"33\u00f76=5, 3\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  // Replace only the first occurrence (results are returned in document\n  // order), mirroring how each source value maps 1:1 onto the next cell\n  // in the table (top-to-bottom, left-to-right).\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# This script applies the set of literal text replacements that correspond\n# to the day's updated date heading and the newly generated division\n# answers in the table. Every run of text is replaced exactly once, using\n# Word's Find/Replace (wdReplaceOne) against the whole document content.\n\n$replacements = @(\n    ,@(\"2025-09-11 Thursday\", \"2025-09-12 Friday\")\n    ,@(\"79\u00f78=9, 7\", \"94\u00f79=10, 4\")\n    ,@(\"49\u00f73=16, 1\", \"74\u00f74=18, 2\")\n    ,@(\"86\u00f79=9, 5\", \"58\u00f77=8, 2\")\n    ,@(\"20\u00f78=2, 4\", \"75\u00f77=10, 5\")\n    ,@(\"52\u00f77=7, 3\", \"70\u00f73=23, 1\")\n    ,@(\"65\u00f79=7, 2\", \"51\u00f74=12, 3\")\n    ,@(\"41\u00f76=6, 5\", \"76\u00f73=25, 1\")\n    ,@(\"66\u00f74=16, 2\", \"22\u00f77=3, 1\")\n    ,@(\"92\u00f73=30, 2\", \"41\u00f76=6, 5\")\n    ,@(\"33\u00f77=4, 5\", \"87\u00f75=17, 2\")\n    ,@(\"75\u00f79=8, 3\", \"43\u00f74=10, 3\")\n    ,@(\"68\u00f74=17, 0\", \"32\u00f73=10, 2\")\n    ,@(\"27\u00f77=3, 6\", \"22\u00f76=3, 4\")\n    ,@(\"15\u00f76=2, 3\", \"16\u00f74=4, 0\")\n    ,@(\"86\u00f72=43, 0\", \"43\u00f77=6, 1\")\n    ,@(\"82\u00f73=27, 1\", \"33\u00f74=8, 1\")\n    ,@(\"65\u00f72=32, 1\", \"23\u00f77=3, 2\")\n    ,@(\"90\u00f75=18, 0\", \"89\u00f74=22, 1\")\n    ,@(\"29\u00f76=4, 5\", \"43\u00f72=21, 1\")\n    ,@(\"34\u00f72=17, 0\", \"65\u00f74=16, 1\")\n    ,@(\"74\u00f73=24, 2\", \"77\u00f75=15, 2\")\n    ,@(\"80\u00f72=40, 0\", \"16\u00f75=3, 1\")\n    ,@(\"13\u00f75=2, 3\", \"31\u00f78=3, 7\")\n    ,@(\"26\u00f77=3, 5\", \"31\u00f72=15, 1\")\n    ,@(\"96\u00f77=13, 5\", \"33\u00f76=5, 3\")\n)\n\n$d = $word.ActiveDocument\n$find = $d.Content.Find\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    # wdFindContinue = 1 (Wrap), wdReplaceOne = 1 (Replace)\n    $found = $find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 1)\n\n    if (-not $found) {\n        throw \"Could not find text to replace: $oldText\"\n    }\n}\n"}
